$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2 qubits")
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 99.81
$ws.Range("D2").Value = 99.66
$ws.Range("E2").Value = 99.48999999999999
$ws.Range("F2").Value = 99.3
$ws.Range("G2").Value = 99.11
$ws.Range("H2").Value = 98.98
$ws.Range("I2").Value = 98.95
$ws.Range("J2").Value = 98.73999999999999
$ws.Range("K2").Value = 98.67
$ws.Range("L2").Value = 98.3
$ws.Range("B3").Value = 50.05
$ws.Range("C3").Value = 50.2
$ws.Range("D3").Value = 48.72
$ws.Range("E3").Value = 49.64
$ws.Range("F3").Value = 49.99
$ws.Range("G3").Value = 49.73
$ws.Range("H3").Value = 49.74
$ws.Range("I3").Value = 50.31
$ws.Range("J3").Value = 49.74
$ws.Range("K3").Value = 50.53
$ws.Range("L3").Value = 50.370000000000005

$ws = $wb.Worksheets.Item("2qbit_oracle")
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2.003807233744114
$ws.Range("D2").Value = 2.006823198876179
$ws.Range("E2").Value = 2.010252286661976
$ws.Range("F2").Value = 2.014098690835851
$ws.Range("G2").Value = 2.017959842599132
$ws.Range("H2").Value = 2.020610224287735
$ws.Range("I2").Value = 2.02122283981809
$ws.Range("J2").Value = 2.02552157180474
$ws.Range("K2").Value = 2.02695854869768
$ws.Range("L2").Value = 2.034587995930824
$ws.Range("B3").Value = 3.996003996003996
$ws.Range("C3").Value = 3.98406374501992
$ws.Range("D3").Value = 4.105090311986864
$ws.Range("E3").Value = 4.029008863819501
$ws.Range("F3").Value = 4.000800160032006
$ws.Range("G3").Value = 4.021717273275689
$ws.Range("H3").Value = 4.020908725371934
$ws.Range("I3").Value = 3.975352812562115
$ws.Range("J3").Value = 4.020908725371934
$ws.Range("K3").Value = 3.958044725905403
$ws.Range("L3").Value = 3.9706174310105222

$ws = $wb.Worksheets.Item("3 qubits")
$ws.Range("B2").Value = 94.64
$ws.Range("C2").Value = 92.72
$ws.Range("D2").Value = 90.88
$ws.Range("E2").Value = 89.22
$ws.Range("F2").Value = 87.37
$ws.Range("G2").Value = 86.16
$ws.Range("H2").Value = 84.06
$ws.Range("I2").Value = 83.2
$ws.Range("J2").Value = 80.89
$ws.Range("K2").Value = 80.16
$ws.Range("L2").Value = 78.62
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 98.87
$ws.Range("D3").Value = 97.87
$ws.Range("E3").Value = 96.94
$ws.Range("F3").Value = 96.22
$ws.Range("G3").Value = 94.81999999999999
$ws.Range("H3").Value = 94.81999999999999
$ws.Range("I3").Value = 92.48999999999999
$ws.Range("J3").Value = 92.62
$ws.Range("K3").Value = 92.62
$ws.Range("L3").Value = 89.64999999999999
$ws.Range("B4").Value = 84.46000000000001
$ws.Range("C4").Value = 84.27
$ws.Range("D4").Value = 83.3
$ws.Range("E4").Value = 82.56999999999999
$ws.Range("F4").Value = 81.82000000000001
$ws.Range("G4").Value = 81.34
$ws.Range("H4").Value = 80.98999999999999
$ws.Range("I4").Value = 80.25999999999999
$ws.Range("J4").Value = 79.67
$ws.Range("K4").Value = 79.44
$ws.Range("L4").Value = 78.75
$ws.Range("B5").Value = 50.03999999999999
$ws.Range("C5").Value = 49.87
$ws.Range("D5").Value = 49.01
$ws.Range("E5").Value = 49.48
$ws.Range("F5").Value = 49.04
$ws.Range("G5").Value = 49.62
$ws.Range("H5").Value = 50.45
$ws.Range("I5").Value = 50.37
$ws.Range("J5").Value = 49.59
$ws.Range("K5").Value = 50.3
$ws.Range("L5").Value = 50.14999999999999

$ws = $wb.Worksheets.Item("3qbit_oracle")
$ws.Range("B2").Value = 3.169907016060862
$ws.Range("C2").Value = 3.235547886108714
$ws.Range("D2").Value = 3.301056338028169
$ws.Range("E2").Value = 3.362474781439139
$ws.Range("F2").Value = 3.433672885429781
$ws.Range("G2").Value = 3.481894150417828
$ws.Range("H2").Value = 3.568879371877231
$ws.Range("I2").Value = 3.605769230769231
$ws.Range("J2").Value = 3.708740264556805
$ws.Range("K2").Value = 3.74251497005988
$ws.Range("L2").Value = 3.815822945815314
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2.022858298776171
$ws.Range("D3").Value = 2.043527127822622
$ws.Range("E3").Value = 2.063131834124201
$ws.Range("F3").Value = 2.078569943878612
$ws.Range("G3").Value = 2.109259649862898
$ws.Range("H3").Value = 2.109259649862898
$ws.Range("I3").Value = 2.162395934695643
$ws.Range("J3").Value = 2.159360829194558
$ws.Range("K3").Value = 2.159360829194558
$ws.Range("L3").Value = 2.230897936419409
$ws.Range("B4").Value = 2.367984844896993
$ws.Range("C4").Value = 2.373323840037973
$ws.Range("D4").Value = 2.400960384153661
$ws.Range("E4").Value = 2.422187235073271
$ws.Range("F4").Value = 2.444390124663896
$ws.Range("G4").Value = 2.458814851241701
$ws.Range("H4").Value = 2.469440671687863
$ws.Range("I4").Value = 2.4919013207077
$ws.Range("J4").Value = 2.51035521526296
$ws.Range("K4").Value = 2.517623363544814
$ws.Range("L4").Value = 2.5396825396825395
$ws.Range("B5").Value = 3.996802557953637
$ws.Range("C5").Value = 4.010427110487267
$ws.Range("D5").Value = 4.080799836768007
$ws.Range("E5").Value = 4.042037186742118
$ws.Range("F5").Value = 4.078303425774878
$ws.Range("G5").Value = 4.030632809351068
$ws.Range("H5").Value = 3.964321110009911
$ws.Range("I5").Value = 3.970617431010522
$ws.Range("J5").Value = 4.033071183706393
$ws.Range("K5").Value = 3.976143141153082
$ws.Range("L5").Value = 3.988035892323031

$ws = $wb.Worksheets.Item("4 qubits")
$ws.Range("B2").Value = 95.95
$ws.Range("C2").Value = 83.11
$ws.Range("D2").Value = 72.13000000000001
$ws.Range("E2").Value = 63.19
$ws.Range("F2").Value = 55.12
$ws.Range("G2").Value = 48.32
$ws.Range("H2").Value = 43.61
$ws.Range("I2").Value = 35.8
$ws.Range("J2").Value = 36.59
$ws.Range("K2").Value = 30.32
$ws.Range("L2").Value = 28.48
$ws.Range("B3").Value = 94.55
$ws.Range("C3").Value = 84.8
$ws.Range("D3").Value = 77.84999999999999
$ws.Range("E3").Value = 71.31999999999999
$ws.Range("F3").Value = 65.91
$ws.Range("G3").Value = 61.8
$ws.Range("H3").Value = 54.83
$ws.Range("I3").Value = 51.39
$ws.Range("J3").Value = 46.97
$ws.Range("K3").Value = 39.48
$ws.Range("L3").Value = 41.010000000000005
$ws.Range("B4").Value = 95.12
$ws.Range("C4").Value = 89.8
$ws.Range("D4").Value = 85.55
$ws.Range("E4").Value = 79.80000000000001
$ws.Range("F4").Value = 76.02
$ws.Range("G4").Value = 72.59999999999999
$ws.Range("H4").Value = 69.98999999999999
$ws.Range("I4").Value = 66.64
$ws.Range("J4").Value = 64.41
$ws.Range("K4").Value = 59.08
$ws.Range("L4").Value = 55.16
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 95.27
$ws.Range("D5").Value = 90.23
$ws.Range("E5").Value = 85.94000000000001
$ws.Range("F5").Value = 82.26000000000001
$ws.Range("G5").Value = 77.84999999999999
$ws.Range("H5").Value = 73.86
$ws.Range("I5").Value = 71.2
$ws.Range("J5").Value = 66.22
$ws.Range("K5").Value = 63.23
$ws.Range("L5").Value = 59.940000000000005
$ws.Range("B6").Value = 95.72
$ws.Range("C6").Value = 91.31
$ws.Range("D6").Value = 87.02
$ws.Range("E6").Value = 82.53
$ws.Range("F6").Value = 79.25999999999999
$ws.Range("G6").Value = 75.70999999999999
$ws.Range("H6").Value = 71.95
$ws.Range("I6").Value = 70.89
$ws.Range("J6").Value = 66.25
$ws.Range("K6").Value = 64.83
$ws.Range("L6").Value = 61.370000000000005
$ws.Range("B7").Value = 84.3
$ws.Range("C7").Value = 80.30000000000001
$ws.Range("D7").Value = 77.94
$ws.Range("E7").Value = 75.94999999999999
$ws.Range("F7").Value = 72.39999999999999
$ws.Range("G7").Value = 70.3
$ws.Range("H7").Value = 67.97999999999999
$ws.Range("I7").Value = 66.16
$ws.Range("J7").Value = 64.68000000000001
$ws.Range("K7").Value = 61.75000000000001
$ws.Range("L7").Value = 60.14000000000001
$ws.Range("B8").Value = 68.15000000000001
$ws.Range("C8").Value = 66.71000000000001
$ws.Range("D8").Value = 64.3
$ws.Range("E8").Value = 63.08000000000001
$ws.Range("F8").Value = 62.82
$ws.Range("G8").Value = 60.79
$ws.Range("H8").Value = 59.92
$ws.Range("I8").Value = 58.8
$ws.Range("J8").Value = 58.15
$ws.Range("K8").Value = 57.17
$ws.Range("L8").Value = 55.669999999999995
$ws.Range("B9").Value = 50.56
$ws.Range("C9").Value = 49.87
$ws.Range("D9").Value = 49.85
$ws.Range("E9").Value = 49.24
$ws.Range("F9").Value = 50.33
$ws.Range("G9").Value = 50.49
$ws.Range("H9").Value = 50.2
$ws.Range("I9").Value = 49.85
$ws.Range("J9").Value = 49.29
$ws.Range("K9").Value = 49.53
$ws.Range("L9").Value = 49.99

$ws = $wb.Worksheets.Item("4qbit_oracle")
$ws.Range("B2").Value = 4.168837936425222
$ws.Range("C2").Value = 4.812898568162676
$ws.Range("D2").Value = 5.545542769998614
$ws.Range("E2").Value = 6.330115524608324
$ws.Range("F2").Value = 7.256894049346879
$ws.Range("G2").Value = 8.278145695364238
$ws.Range("H2").Value = 9.172208209126348
$ws.Range("I2").Value = 11.1731843575419
$ws.Range("J2").Value = 10.93194861984149
$ws.Range("K2").Value = 13.19261213720317
$ws.Range("L2").Value = 14.044943820224718
$ws.Range("B3").Value = 3.172924378635642
$ws.Range("C3").Value = 3.537735849056604
$ws.Range("D3").Value = 3.853564547206166
$ws.Range("E3").Value = 4.206393718452048
$ws.Range("F3").Value = 4.551661356395084
$ws.Range("G3").Value = 4.854368932038836
$ws.Range("H3").Value = 5.47145723144264
$ws.Range("I3").Value = 5.837711617046118
$ws.Range("J3").Value = 6.387055567383436
$ws.Range("K3").Value = 7.598784194528876
$ws.Range("L3").Value = 7.315288953913679
$ws.Range("B4").Value = 2.102607232968881
$ws.Range("C4").Value = 2.2271714922049
$ws.Range("D4").Value = 2.33781414377557
$ws.Range("E4").Value = 2.506265664160401
$ws.Range("F4").Value = 2.630886608787161
$ws.Range("G4").Value = 2.754820936639118
$ws.Range("H4").Value = 2.857551078725532
$ws.Range("I4").Value = 3.001200480192077
$ws.Range("J4").Value = 3.105107902499612
$ws.Range("K4").Value = 3.385240352064997
$ws.Range("L4").Value = 3.625815808556925
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2.099296735593576
$ws.Range("D5").Value = 2.216557685913776
$ws.Range("E5").Value = 2.327205026762858
$ws.Range("F5").Value = 2.431315341599805
$ws.Range("G5").Value = 2.569043031470777
$ws.Range("H5").Value = 2.707825616030328
$ws.Range("I5").Value = 2.808988764044944
$ws.Range("J5").Value = 3.020235578375113
$ws.Range("K5").Value = 3.163055511624229
$ws.Range("L5").Value = 3.33667000333667
$ws.Range("B6").Value = 2.089427496865859
$ws.Range("C6").Value = 2.190340597962983
$ws.Range("D6").Value = 2.298322224775914
$ws.Range("E6").Value = 2.423361201987156
$ws.Range("F6").Value = 2.523340903356043
$ws.Range("G6").Value = 2.641658961828028
$ws.Range("H6").Value = 2.779708130646282
$ws.Range("I6").Value = 2.821272393849626
$ws.Range("J6").Value = 3.018867924528302
$ws.Range("K6").Value = 3.08499151627333
$ws.Range("L6").Value = 3.258921297050676
$ws.Range("B7").Value = 2.372479240806643
$ws.Range("C7").Value = 2.4906600249066
$ws.Range("D7").Value = 2.566076469078779
$ws.Range("E7").Value = 2.633311389071758
$ws.Range("F7").Value = 2.762430939226519
$ws.Range("G7").Value = 2.844950213371266
$ws.Range("H7").Value = 2.942041776993233
$ws.Range("I7").Value = 3.022974607013301
$ws.Range("J7").Value = 3.092145949288807
$ws.Range("K7").Value = 3.238866396761133
$ws.Range("L7").Value = 3.3255736614566014
$ws.Range("B8").Value = 2.93470286133529
$ws.Range("C8").Value = 2.99805126667666
$ws.Range("D8").Value = 3.110419906687403
$ws.Range("E8").Value = 3.170577045022194
$ws.Range("F8").Value = 3.183699458771092
$ws.Range("G8").Value = 3.290014805066623
$ws.Range("H8").Value = 3.337783711615488
$ws.Range("I8").Value = 3.401360544217687
$ws.Range("J8").Value = 3.439380911435942
$ws.Range("K8").Value = 3.498338289312577
$ws.Range("L8").Value = 3.592599245554158
$ws.Range("B9").Value = 3.955696202531646
$ws.Range("C9").Value = 4.010427110487267
$ws.Range("D9").Value = 4.012036108324975
$ws.Range("E9").Value = 4.061738424045491
$ws.Range("F9").Value = 3.973773097556129
$ws.Range("G9").Value = 3.961180431768667
$ws.Range("H9").Value = 3.98406374501992
$ws.Range("I9").Value = 4.012036108324975
$ws.Range("J9").Value = 4.057618178129438
$ws.Range("K9").Value = 4.037956793862306
$ws.Range("L9").Value = 4.000800160032006

$ws = $wb.Worksheets.Item("5 qubits")
$ws.Range("B2").Value = 99.92
$ws.Range("C2").Value = 43.31
$ws.Range("D2").Value = 19.69
$ws.Range("E2").Value = 11.39
$ws.Range("F2").Value = 8.129999999999999
$ws.Range("G2").Value = 4.82
$ws.Range("H2").Value = 4.79
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 3.96
$ws.Range("K2").Value = 3.51
$ws.Range("L2").Value = 3.27
$ws.Range("B3").Value = 96.19
$ws.Range("C3").Value = 49.71
$ws.Range("D3").Value = 26.42
$ws.Range("E3").Value = 17.51
$ws.Range("F3").Value = 11.25
$ws.Range("G3").Value = 8.49
$ws.Range("H3").Value = 8.41
$ws.Range("I3").Value = 6.98
$ws.Range("J3").Value = 7.06
$ws.Range("K3").Value = 6.859999999999999
$ws.Range("L3").Value = 6.460000000000001
$ws.Range("B4").Value = 99.95
$ws.Range("C4").Value = 61.05
$ws.Range("D4").Value = 38.34
$ws.Range("E4").Value = 29.18
$ws.Range("F4").Value = 21.99
$ws.Range("G4").Value = 15.11
$ws.Range("H4").Value = 14.87
$ws.Range("I4").Value = 12.65
$ws.Range("J4").Value = 10.98
$ws.Range("K4").Value = 10.87
$ws.Range("L4").Value = 9.84
$ws.Range("B5").Value = 94.73
$ws.Range("C5").Value = 57.4
$ws.Range("D5").Value = 40.58
$ws.Range("E5").Value = 27.8
$ws.Range("F5").Value = 20.58
$ws.Range("G5").Value = 19.02
$ws.Range("H5").Value = 15.06
$ws.Range("I5").Value = 14.69
$ws.Range("J5").Value = 14.78
$ws.Range("K5").Value = 13.51
$ws.Range("L5").Value = 13.020000000000001
$ws.Range("B6").Value = 87.78
$ws.Range("C6").Value = 67.95
$ws.Range("D6").Value = 53.02
$ws.Range("E6").Value = 43.62
$ws.Range("F6").Value = 33.66
$ws.Range("G6").Value = 30.09
$ws.Range("H6").Value = 26.14
$ws.Range("I6").Value = 23.47
$ws.Range("J6").Value = 21.04
$ws.Range("K6").Value = 19.39
$ws.Range("L6").Value = 18.86
$ws.Range("B7").Value = 94.63
$ws.Range("C7").Value = 74.11999999999999
$ws.Range("D7").Value = 58.22000000000001
$ws.Range("E7").Value = 47.09
$ws.Range("F7").Value = 38.44
$ws.Range("G7").Value = 33.77
$ws.Range("H7").Value = 29.82
$ws.Range("I7").Value = 26.23
$ws.Range("J7").Value = 24.03
$ws.Range("K7").Value = 23.52
$ws.Range("L7").Value = 22.13
$ws.Range("B8").Value = 98.87
$ws.Range("C8").Value = 77.45999999999999
$ws.Range("D8").Value = 60.89
$ws.Range("E8").Value = 52.58000000000001
$ws.Range("F8").Value = 41.87
$ws.Range("G8").Value = 36.01
$ws.Range("H8").Value = 31.53
$ws.Range("I8").Value = 29.68
$ws.Range("J8").Value = 27.25
$ws.Range("K8").Value = 26.16
$ws.Range("L8").Value = 26.26
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 78.51000000000001
$ws.Range("D9").Value = 63.54
$ws.Range("E9").Value = 52.11
$ws.Range("F9").Value = 44.52
$ws.Range("G9").Value = 39.15
$ws.Range("H9").Value = 35.70999999999999
$ws.Range("I9").Value = 32.97
$ws.Range("J9").Value = 31.56
$ws.Range("K9").Value = 29.41
$ws.Range("L9").Value = 27.500000000000004
$ws.Range("B10").Value = 98.91
$ws.Range("C10").Value = 78.45
$ws.Range("D10").Value = 63.36000000000001
$ws.Range("E10").Value = 53.22
$ws.Range("F10").Value = 47.03
$ws.Range("G10").Value = 41.63
$ws.Range("H10").Value = 37.69
$ws.Range("I10").Value = 34.83
$ws.Range("J10").Value = 33.1
$ws.Range("K10").Value = 31.25
$ws.Range("L10").Value = 30.220000000000002
$ws.Range("B11").Value = 95.91
$ws.Range("C11").Value = 76.7
$ws.Range("D11").Value = 64.07000000000001
$ws.Range("E11").Value = 54.26
$ws.Range("F11").Value = 47.7
$ws.Range("G11").Value = 44.83
$ws.Range("H11").Value = 39.68
$ws.Range("I11").Value = 37.08
$ws.Range("J11").Value = 35.3
$ws.Range("K11").Value = 34.37
$ws.Range("L11").Value = 34.449999999999996
$ws.Range("B12").Value = 91.06
$ws.Range("C12").Value = 74.56
$ws.Range("D12").Value = 63.05
$ws.Range("E12").Value = 54.79000000000001
$ws.Range("F12").Value = 48.52
$ws.Range("G12").Value = 43.82
$ws.Range("H12").Value = 41.3
$ws.Range("I12").Value = 39.06
$ws.Range("J12").Value = 38.82
$ws.Range("K12").Value = 37.27
$ws.Range("L12").Value = 35.6
$ws.Range("B13").Value = 84.66
$ws.Range("C13").Value = 71.86
$ws.Range("D13").Value = 61.17
$ws.Range("E13").Value = 53.63
$ws.Range("F13").Value = 49.17
$ws.Range("G13").Value = 45.72
$ws.Range("H13").Value = 43.64
$ws.Range("I13").Value = 42.64
$ws.Range("J13").Value = 39.84
$ws.Range("K13").Value = 40.28
$ws.Range("L13").Value = 38.769999999999996
$ws.Range("B14").Value = 76.99000000000001
$ws.Range("C14").Value = 66.19
$ws.Range("D14").Value = 58.54000000000001
$ws.Range("E14").Value = 52.87
$ws.Range("F14").Value = 49.64
$ws.Range("G14").Value = 47.29
$ws.Range("H14").Value = 45.09
$ws.Range("I14").Value = 43.59
$ws.Range("J14").Value = 43.82
$ws.Range("K14").Value = 42.38
$ws.Range("L14").Value = 42.08
$ws.Range("B15").Value = 68.17999999999999
$ws.Range("C15").Value = 61.49
$ws.Range("D15").Value = 56.14
$ws.Range("E15").Value = 52.43
$ws.Range("F15").Value = 50.45
$ws.Range("G15").Value = 48.49
$ws.Range("H15").Value = 46.86
$ws.Range("I15").Value = 46.43
$ws.Range("J15").Value = 44.59
$ws.Range("K15").Value = 45.18
$ws.Range("L15").Value = 44.76
$ws.Range("B16").Value = 59.21999999999999
$ws.Range("C16").Value = 56.48
$ws.Range("D16").Value = 53.14
$ws.Range("E16").Value = 50.71
$ws.Range("F16").Value = 50.16
$ws.Range("G16").Value = 49.22000000000001
$ws.Range("H16").Value = 48.19
$ws.Range("I16").Value = 47.73
$ws.Range("J16").Value = 47.54
$ws.Range("K16").Value = 47.77
$ws.Range("L16").Value = 47.839999999999996
$ws.Range("B17").Value = 50.62
$ws.Range("C17").Value = 49.68
$ws.Range("D17").Value = 50.26000000000001
$ws.Range("E17").Value = 49.81
$ws.Range("F17").Value = 49.15
$ws.Range("G17").Value = 48.93
$ws.Range("H17").Value = 50.29
$ws.Range("I17").Value = 49.08
$ws.Range("J17").Value = 49.5
$ws.Range("K17").Value = 49.28
$ws.Range("L17").Value = 50.839999999999996

$ws = $wb.Worksheets.Item("5qbit_oracle")
$ws.Range("B2").Value = 5.00400320256205
$ws.Range("C2").Value = 11.54467790348649
$ws.Range("D2").Value = 25.39360081259522
$ws.Range("E2").Value = 43.89815627743634
$ws.Range("F2").Value = 61.50061500615006
$ws.Range("G2").Value = 103.7344398340249
$ws.Range("H2").Value = 104.384133611691
$ws.Range("I2").Value = 138.8888888888889
$ws.Range("J2").Value = 126.2626262626263
$ws.Range("K2").Value = 142.4501424501424
$ws.Range("L2").Value = 152.90519877675843
$ws.Range("B3").Value = 4.158436427903109
$ws.Range("C3").Value = 8.046670690002012
$ws.Range("D3").Value = 15.14004542013626
$ws.Range("E3").Value = 22.84408909194746
$ws.Range("F3").Value = 35.55555555555556
$ws.Range("G3").Value = 47.11425206124853
$ws.Range("H3").Value = 47.56242568370987
$ws.Range("I3").Value = 57.30659025787966
$ws.Range("J3").Value = 56.657223796034
$ws.Range("K3").Value = 58.30903790087464
$ws.Range("L3").Value = 61.919504643962846
$ws.Range("B4").Value = 3.001500750375188
$ws.Range("C4").Value = 4.914004914004914
$ws.Range("D4").Value = 7.82472613458529
$ws.Range("E4").Value = 10.28101439342015
$ws.Range("F4").Value = 13.64256480218281
$ws.Range("G4").Value = 19.85440105890139
$ws.Range("H4").Value = 20.17484868863484
$ws.Range("I4").Value = 23.71541501976284
$ws.Range("J4").Value = 27.3224043715847
$ws.Range("K4").Value = 27.59889604415823
$ws.Range("L4").Value = 30.487804878048777
$ws.Range("B5").Value = 3.166895386889053
$ws.Range("C5").Value = 5.226480836236933
$ws.Range("D5").Value = 7.392804337111878
$ws.Range("E5").Value = 10.79136690647482
$ws.Range("F5").Value = 14.57725947521866
$ws.Range("G5").Value = 15.77287066246057
$ws.Range("H5").Value = 19.9203187250996
$ws.Range("I5").Value = 20.42205582028591
$ws.Range("J5").Value = 20.29769959404601
$ws.Range("K5").Value = 22.20577350111029
$ws.Range("L5").Value = 23.04147465437788
$ws.Range("B6").Value = 2.27842333105491
$ws.Range("C6").Value = 2.943340691685062
$ws.Range("D6").Value = 3.772161448509997
$ws.Range("E6").Value = 4.585052728106374
$ws.Range("F6").Value = 5.941770647653001
$ws.Range("G6").Value = 6.646726487205052
$ws.Range("H6").Value = 7.651109410864575
$ws.Range("I6").Value = 8.52151682999574
$ws.Range("J6").Value = 9.505703422053232
$ws.Range("K6").Value = 10.31459515214028
$ws.Range("L6").Value = 10.604453870625663
$ws.Range("B7").Value = 2.113494663425975
$ws.Range("C7").Value = 2.698327037236913
$ws.Range("D7").Value = 3.435245620061834
$ws.Range("E7").Value = 4.247186239116585
$ws.Range("F7").Value = 5.202913631633715
$ws.Range("G7").Value = 5.922416345869115
$ws.Range("H7").Value = 6.70690811535882
$ws.Range("I7").Value = 7.624857033930613
$ws.Range("J7").Value = 8.322929671244278
$ws.Range("K7").Value = 8.503401360544217
$ws.Range("L7").Value = 9.03750564844103
$ws.Range("B8").Value = 2.022858298776171
$ws.Range("C8").Value = 2.581977794990963
$ws.Range("D8").Value = 3.284611594678929
$ws.Range("E8").Value = 3.803727653100038
$ws.Range("F8").Value = 4.776689754000477
$ws.Range("G8").Value = 5.55401277422938
$ws.Range("H8").Value = 6.343165239454488
$ws.Range("I8").Value = 6.738544474393531
$ws.Range("J8").Value = 7.339449541284404
$ws.Range("K8").Value = 7.64525993883792
$ws.Range("L8").Value = 7.616146230007616
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2.547446185199337
$ws.Range("D9").Value = 3.147623544224111
$ws.Range("E9").Value = 3.838034926117828
$ws.Range("F9").Value = 4.492362982929021
$ws.Range("G9").Value = 5.108556832694764
$ws.Range("H9").Value = 5.600672080649678
$ws.Range("I9").Value = 6.066120715802245
$ws.Range("J9").Value = 6.337135614702155
$ws.Range("K9").Value = 6.800408024481469
$ws.Range("L9").Value = 7.2727272727272725
$ws.Range("B10").Value = 2.022040238600748
$ws.Range("C10").Value = 2.549394518801785
$ws.Range("D10").Value = 3.156565656565657
$ws.Range("E10").Value = 3.757985719654265
$ws.Range("F10").Value = 4.25260472039124
$ws.Range("G10").Value = 4.804227720393946
$ws.Range("H10").Value = 5.306447333510215
$ws.Range("I10").Value = 5.742176284811944
$ws.Range("J10").Value = 6.042296072507553
$ws.Range("K10").Value = 6.4
$ws.Range("L10").Value = 6.618133686300463
$ws.Range("B11").Value = 2.085288291106246
$ws.Range("C11").Value = 2.607561929595828
$ws.Range("D11").Value = 3.121585765568909
$ws.Range("E11").Value = 3.685956505713233
$ws.Range("F11").Value = 4.19287211740042
$ws.Range("G11").Value = 4.461298237787196
$ws.Range("H11").Value = 5.040322580645161
$ws.Range("I11").Value = 5.393743257820928
$ws.Range("J11").Value = 5.665722379603399
$ws.Range("K11").Value = 5.819028222286878
$ws.Range("L11").Value = 5.805515239477503
$ws.Range("B12").Value = 2.196354052273227
$ws.Range("C12").Value = 2.682403433476394
$ws.Range("D12").Value = 3.17208564631245
$ws.Range("E12").Value = 3.650301149844862
$ws.Range("F12").Value = 4.122011541632316
$ws.Range("G12").Value = 4.564125969876769
$ws.Range("H12").Value = 4.842615012106537
$ws.Range("I12").Value = 5.120327700972863
$ws.Range("J12").Value = 5.151983513652756
$ws.Range("K12").Value = 5.366246310705661
$ws.Range("L12").Value = 5.617977528089888
$ws.Range("B13").Value = 2.362390739428301
$ws.Range("C13").Value = 2.783189535207347
$ws.Range("D13").Value = 3.269576589831617
$ws.Range("E13").Value = 3.729256013425322
$ws.Range("F13").Value = 4.067520846044336
$ws.Range("G13").Value = 4.374453193350831
$ws.Range("H13").Value = 4.582951420714941
$ws.Range("I13").Value = 4.690431519699812
$ws.Range("J13").Value = 5.020080321285141
$ws.Range("K13").Value = 4.965243296921549
$ws.Range("L13").Value = 5.158627805003869
$ws.Range("B14").Value = 2.59773996622938
$ws.Range("C14").Value = 3.021604471974618
$ws.Range("D14").Value = 3.41646737273659
$ws.Range("E14").Value = 3.782863627766219
$ws.Range("F14").Value = 4.029008863819501
$ws.Range("G14").Value = 4.229223937407486
$ws.Range("H14").Value = 4.435573297848747
$ws.Range("I14").Value = 4.588208304657031
$ws.Range("J14").Value = 4.564125969876769
$ws.Range("K14").Value = 4.719207173194904
$ws.Range("L14").Value = 4.752851711026616
$ws.Range("B15").Value = 2.933411557641537
$ws.Range("C15").Value = 3.252561392096276
$ws.Range("D15").Value = 3.562522265764161
$ws.Range("E15").Value = 3.814609956131986
$ws.Range("F15").Value = 3.964321110009911
$ws.Range("G15").Value = 4.124561765312436
$ws.Range("H15").Value = 4.268032437046521
$ws.Range("I15").Value = 4.307559767391773
$ws.Range("J15").Value = 4.485310607759588
$ws.Range("K15").Value = 4.426737494466578
$ws.Range("L15").Value = 4.468275245755138
$ws.Range("B16").Value = 3.377237419790611
$ws.Range("C16").Value = 3.541076487252125
$ws.Range("D16").Value = 3.763643206624012
$ws.Range("E16").Value = 3.94399526720568
$ws.Range("F16").Value = 3.987240829346093
$ws.Range("G16").Value = 4.063388866314506
$ws.Range("H16").Value = 4.150238638721727
$ws.Range("I16").Value = 4.190236748376283
$ws.Range("J16").Value = 4.206983592763988
$ws.Range("K16").Value = 4.186728072011723
$ws.Range("L16").Value = 4.1806020066889635
$ws.Range("B17").Value = 3.951007506914263
$ws.Range("C17").Value = 4.025764895330113
$ws.Range("D17").Value = 3.979307600477517
$ws.Range("E17").Value = 4.015257980325236
$ws.Range("F17").Value = 4.069175991861648
$ws.Range("G17").Value = 4.087471898630697
$ws.Range("H17").Value = 3.976933784052496
$ws.Range("I17").Value = 4.074979625101874
$ws.Range("J17").Value = 4.040404040404041
$ws.Range("K17").Value = 4.058441558441558
$ws.Range("L17").Value = 3.933910306845004
